# Add two new rows ("study-summary" and "study-variable-summary" profiles)
# to the Observations worksheet, mirroring the formatting of the existing
# data row (row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3: study-summary ----
$ws.Range("A2:K2").Copy()
$ws.Range("A3:K3").PasteSpecial(-4122)   # xlPasteFormats - carry over the row's style

$ws.Range("A3").Value2 = "study-summary"
$ws.Range("B3").Value2 = "Study Summary"
$ws.Range("E3").Value2 = "null#C0242482"
$ws.Range("F3").Value2 = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G3").Value2 = "dateTime, Period, Timing, instant"
$ws.Range("H3").Value2 = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I3").Value2 = "optional"

# ---- Row 4: study-variable-summary ----
$ws.Range("A2:K2").Copy()
$ws.Range("A4:K4").PasteSpecial(-4122)   # xlPasteFormats - carry over the row's style

$ws.Range("A4").Value2 = "study-variable-summary"
$ws.Range("B4").Value2 = "Study Variable Summary"
$ws.Range("E4").Value2 = "null#C0242482"
$ws.Range("F4").Value2 = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G4").Value2 = "dateTime, Period, Timing, instant"
$ws.Range("H4").Value2 = "CodeableConcept"
$ws.Range("I4").Value2 = "optional"
